# Fruta / hortaliza, semanal
#
# A new weekly record was inserted into the "Ají" (pepper) sheet as row 347,
# pushing the existing rows 347-378 down to 348-379 (dimension grows from
# A1:R378 to A1:R379).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 347 - shifts old rows 347:378 down to 348:379
# and Excel extends the used range / dimension automatically.
$ws.Rows.Item(347).Insert()

# Populate the new row with the latest weekly price record.
$ws.Cells.Item(347, 1).Value  = 5
$ws.Cells.Item(347, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(347, 3).Value  = "Maule"
$ws.Cells.Item(347, 4).Value  = 45132
$ws.Cells.Item(347, 5).Value  = 7
$ws.Cells.Item(347, 6).Value  = 100112021
$ws.Cells.Item(347, 7).Value  = "Ají"
$ws.Cells.Item(347, 8).Value  = "Inferno"
$ws.Cells.Item(347, 9).Value  = "Primera"
$ws.Cells.Item(347, 10).Value = 200
$ws.Cells.Item(347, 11).Value = 12000
$ws.Cells.Item(347, 12).Value = 12000
$ws.Cells.Item(347, 13).Value = 12000
$ws.Cells.Item(347, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(347, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(347, 16).Value = 1200
$ws.Cells.Item(347, 17).Value = 10
$ws.Cells.Item(347, 18).Value = "Hortaliza"
